$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# A new handoff cycle completed for "9a6ee2a4-ee99-45d4-8aa5-40abe742a38f.md".
# The status-report rows get re-emitted in a new order (the freshly
# regenerated row for 9a6ee2a4 moves to the end of each table) and every
# row's Status flips from "Handed back: in sync with en-US" to
# "Ready for handoff". This script rewrites every affected cell explicitly.
# ---------------------------------------------------------------------------

$longError = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/90b130cb884d132a14ef9b598a691df82ff3f386/e2e/9a6ee2a4-ee99-45d4-8aa5-40abe742a38f.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b69e20bc3173fb32d604e623b149cd17b33c27f2/e2e/9a6ee2a4-ee99-45d4-8aa5-40abe742a38f.md."

# ===========================================================================
# Sheet "Overview"
# ===========================================================================
$ws = $wb.Worksheets.Item("Overview")

# Row 2 <- old row 3 (ffffae36bc29...)
$ws.Cells.Item(2,1).Value = "ffffae36bc29-3452-4bfd-b460-2df2bf5dd1f4.md"
$ws.Cells.Item(2,2).Value = "e2e\ffffae36bc29-3452-4bfd-b460-2df2bf5dd1f4.md"
$ws.Cells.Item(2,3).Value = ".md"
$ws.Cells.Item(2,5).Value = "Handed back: in sync with en-US"
$ws.Cells.Item(2,6).Value = "Handed back: in sync with en-US"
$ws.Cells.Item(2,7).Value = "2016-08-22 07:00:59"

# Row 3 <- old row 4 (ffffffa0ee2988...)
$ws.Cells.Item(3,1).Value = "ffffffa0ee2988-5825-4b77-936c-7f2e8b1b5237.md"
$ws.Cells.Item(3,2).Value = "e2e\ffffffa0ee2988-5825-4b77-936c-7f2e8b1b5237.md"
$ws.Cells.Item(3,3).Value = ".md"
$ws.Cells.Item(3,5).Value = "Handed back: in sync with en-US"
$ws.Cells.Item(3,6).Value = "Handed back: in sync with en-US"
$ws.Cells.Item(3,7).Value = "2016-08-22 07:00:59"

# Row 4 <- old row 2 (9a6ee2a4...), refreshed handoff
$ws.Cells.Item(4,1).Value = "9a6ee2a4-ee99-45d4-8aa5-40abe742a38f.md"
$ws.Cells.Item(4,2).Value = "e2e\9a6ee2a4-ee99-45d4-8aa5-40abe742a38f.md"
$ws.Cells.Item(4,3).Value = ".md"
$ws.Cells.Item(4,5).Value = "Ready for handoff"
$ws.Cells.Item(4,6).Value = "Ready for handoff"
$ws.Cells.Item(4,7).Value = "2016-08-22 07:03:23"

foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$2') {
        $hl.TextToDisplay = "e2e\ffffae36bc29-3452-4bfd-b460-2df2bf5dd1f4.md"
    } elseif ($addr -eq '$B$3') {
        $hl.TextToDisplay = "e2e\ffffffa0ee2988-5825-4b77-936c-7f2e8b1b5237.md"
    } elseif ($addr -eq '$B$4') {
        $hl.TextToDisplay = "e2e\9a6ee2a4-ee99-45d4-8aa5-40abe742a38f.md"
    }
}

# ===========================================================================
# Sheet "zh-cn"
# ===========================================================================
$ws = $wb.Worksheets.Item("zh-cn")

# Row 2 <- old row 3 (ffffae36bc29...)
$ws.Cells.Item(2,1).Value  = "ffffae36bc29-3452-4bfd-b460-2df2bf5dd1f4.md"
$ws.Cells.Item(2,3).Value  = "Ready for handoff"
$ws.Cells.Item(2,6).Value  = "False"
$ws.Cells.Item(2,7).Value  = "17881179-f98c-4b66-b862-2e75806c0854.1319fb9b3412fc569f2cef548eed7cec62d0c61e.zh-cn.xlf"
$ws.Cells.Item(2,8).Value  = "2016-08-22 07:00:52"
$ws.Cells.Item(2,9).Value  = "17881179-f98c-4b66-b862-2e75806c0854.md"
$ws.Cells.Item(2,10).Value = "17881179-f98c-4b66-b862-2e75806c0854.1319fb9b3412fc569f2cef548eed7cec62d0c61e.zh-cn.xlf"
$ws.Cells.Item(2,11).Value = "2016-08-22 07:01:26"

# Row 3 <- old row 4 (ffffffa0ee2988...)
$ws.Cells.Item(3,1).Value  = "ffffffa0ee2988-5825-4b77-936c-7f2e8b1b5237.md"
$ws.Cells.Item(3,3).Value  = "Ready for handoff"
$ws.Cells.Item(3,6).Value  = "True"
$ws.Cells.Item(3,7).Value  = "17881179-f98c-4b66-b862-2e75806c0854.1319fb9b3412fc569f2cef548eed7cec62d0c61e.zh-cn.xlf"
$ws.Cells.Item(3,8).Value  = "2016-08-22 07:00:52"
$ws.Cells.Item(3,9).Value  = "17881179-f98c-4b66-b862-2e75806c0854.md"
$ws.Cells.Item(3,10).Value = "17881179-f98c-4b66-b862-2e75806c0854.1319fb9b3412fc569f2cef548eed7cec62d0c61e.zh-cn.xlf"
$ws.Cells.Item(3,11).Value = "2016-08-22 07:01:26"

# Row 4 <- old row 2 (9a6ee2a4...), refreshed handoff
$ws.Cells.Item(4,1).Value  = "9a6ee2a4-ee99-45d4-8aa5-40abe742a38f.md"
$ws.Cells.Item(4,3).Value  = "Ready for handoff"
$ws.Cells.Item(4,6).Value  = "False"
$ws.Cells.Item(4,7).Value  = "9a6ee2a4-ee99-45d4-8aa5-40abe742a38f.c338a5143f5511ac63b236d515d6010130c78c36.zh-cn.xlf"
$ws.Cells.Item(4,8).Value  = "2016-08-22 07:03:19"
$ws.Cells.Item(4,9).Value  = "9a6ee2a4-ee99-45d4-8aa5-40abe742a38f.md"
$ws.Cells.Item(4,10).Value = "9a6ee2a4-ee99-45d4-8aa5-40abe742a38f.c338a5143f5511ac63b236d515d6010130c78c36.zh-cn.xlf"
$ws.Cells.Item(4,11).Value = "2016-08-22 07:02:48"
$ws.Cells.Item(4,16).Value = $longError

foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "ffffae36bc29-3452-4bfd-b460-2df2bf5dd1f4.md"
    } elseif ($addr -eq '$I$2') {
        $hl.TextToDisplay = "17881179-f98c-4b66-b862-2e75806c0854.md"
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = "ffffffa0ee2988-5825-4b77-936c-7f2e8b1b5237.md"
    } elseif ($addr -eq '$I$3') {
        $hl.TextToDisplay = "17881179-f98c-4b66-b862-2e75806c0854.md"
    } elseif ($addr -eq '$A$4') {
        $hl.TextToDisplay = "9a6ee2a4-ee99-45d4-8aa5-40abe742a38f.md"
    } elseif ($addr -eq '$I$4') {
        $hl.TextToDisplay = "9a6ee2a4-ee99-45d4-8aa5-40abe742a38f.md"
    }
}

# Column P (16) widened to fit the long error-detail text
$ws.Columns.Item(16).ColumnWidth = 39.17

# ===========================================================================
# Sheet "de-de"
# ===========================================================================
$ws = $wb.Worksheets.Item("de-de")

# Row 2 <- old row 3 (ffffae36bc29...)
$ws.Cells.Item(2,1).Value  = "ffffae36bc29-3452-4bfd-b460-2df2bf5dd1f4.md"
$ws.Cells.Item(2,3).Value  = "Ready for handoff"
$ws.Cells.Item(2,6).Value  = "False"
$ws.Cells.Item(2,7).Value  = "17881179-f98c-4b66-b862-2e75806c0854.1319fb9b3412fc569f2cef548eed7cec62d0c61e.de-de.xlf"
$ws.Cells.Item(2,8).Value  = "2016-08-22 07:00:59"
$ws.Cells.Item(2,9).Value  = "17881179-f98c-4b66-b862-2e75806c0854.md"
$ws.Cells.Item(2,10).Value = "17881179-f98c-4b66-b862-2e75806c0854.1319fb9b3412fc569f2cef548eed7cec62d0c61e.de-de.xlf"
$ws.Cells.Item(2,11).Value = "2016-08-22 07:01:32"

# Row 3 <- old row 4 (ffffffa0ee2988...)
$ws.Cells.Item(3,1).Value  = "ffffffa0ee2988-5825-4b77-936c-7f2e8b1b5237.md"
$ws.Cells.Item(3,3).Value  = "Ready for handoff"
$ws.Cells.Item(3,6).Value  = "True"
$ws.Cells.Item(3,7).Value  = "17881179-f98c-4b66-b862-2e75806c0854.1319fb9b3412fc569f2cef548eed7cec62d0c61e.de-de.xlf"
$ws.Cells.Item(3,8).Value  = "2016-08-22 07:00:59"
$ws.Cells.Item(3,9).Value  = "17881179-f98c-4b66-b862-2e75806c0854.md"
$ws.Cells.Item(3,10).Value = "17881179-f98c-4b66-b862-2e75806c0854.1319fb9b3412fc569f2cef548eed7cec62d0c61e.de-de.xlf"
$ws.Cells.Item(3,11).Value = "2016-08-22 07:01:32"

# Row 4 <- old row 2 (9a6ee2a4...), refreshed handoff
$ws.Cells.Item(4,1).Value  = "9a6ee2a4-ee99-45d4-8aa5-40abe742a38f.md"
$ws.Cells.Item(4,3).Value  = "Ready for handoff"
$ws.Cells.Item(4,6).Value  = "False"
$ws.Cells.Item(4,7).Value  = "9a6ee2a4-ee99-45d4-8aa5-40abe742a38f.c338a5143f5511ac63b236d515d6010130c78c36.de-de.xlf"
$ws.Cells.Item(4,8).Value  = "2016-08-22 07:03:23"
$ws.Cells.Item(4,9).Value  = "9a6ee2a4-ee99-45d4-8aa5-40abe742a38f.md"
$ws.Cells.Item(4,10).Value = "9a6ee2a4-ee99-45d4-8aa5-40abe742a38f.c338a5143f5511ac63b236d515d6010130c78c36.de-de.xlf"
$ws.Cells.Item(4,11).Value = "2016-08-22 07:02:55"
$ws.Cells.Item(4,16).Value = $longError

foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "ffffae36bc29-3452-4bfd-b460-2df2bf5dd1f4.md"
    } elseif ($addr -eq '$I$2') {
        $hl.TextToDisplay = "17881179-f98c-4b66-b862-2e75806c0854.md"
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = "ffffffa0ee2988-5825-4b77-936c-7f2e8b1b5237.md"
    } elseif ($addr -eq '$I$3') {
        $hl.TextToDisplay = "17881179-f98c-4b66-b862-2e75806c0854.md"
    } elseif ($addr -eq '$A$4') {
        $hl.TextToDisplay = "9a6ee2a4-ee99-45d4-8aa5-40abe742a38f.md"
    } elseif ($addr -eq '$I$4') {
        $hl.TextToDisplay = "9a6ee2a4-ee99-45d4-8aa5-40abe742a38f.md"
    }
}

# Column P (16) widened to fit the long error-detail text
$ws.Columns.Item(16).ColumnWidth = 39.17
